# Weekly data refresh: a new observation is prepended to the data block
# (row 32), shifting all existing data rows (old 32..144) down by one
# (new 33..145). This mirrors the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 32, pushing rows 32..144 down to 33..145.
$ws.Rows.Item(32).Insert()

# Populate the new row 32 with this week's record. Columns A,B,C,E,F,G,H,I,
# N,O,Q,R carry the same constants as the rest of the "Orégano" block;
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) are the new values.
$ws.Range("A32").Value = 6
$ws.Range("B32").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 44575
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = 100112029
$ws.Range("G32").Value = "Orégano"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 33
$ws.Range("K32").Value = 8000
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = 8424
$ws.Range("N32").Value = "`$/docena de atados"
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 2808
$ws.Range("Q32").Value = 3
$ws.Range("R32").Value = "Hortaliza"

# Make sure the new D32 cell carries the same date number format as the
# rest of column D (style index 2 in the original file).
$ws.Range("D32").NumberFormat = $ws.Range("D33").NumberFormat
